# chore: added git code in the flow diagram
#
# Applies the OOXML diff to flow.pptx (single-slide excerpt):
#   1. Grows "TextBox 16" (id=17) and appends the new git-push instructions.
#   2. Deletes "TextBox 19" (id=20) which duplicated that same content.
#   3. Bolds the npm command inside "TextBox 20" (id=21) and tweaks spacing.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "TextBox 16" - grow the box and add the new paragraphs describing the
#    git push / VS Code command-palette workflow.
# ---------------------------------------------------------------------------
$tb16 = $s.Shapes.Item(14)

$tf16 = $tb16.TextFrame
$tr16 = $tf16.TextRange

$newParas = "`r" + `
    "`r" + `
    "git push --follow-tags origin master" + `
    "`r" + `
    "`r" + `
    "In VS Code, Push changes with tags by pressing Ctrl + Shift + P" + `
    "`r" + `
    "Type Git: Push (Follow tags)" + `
    "`r" + `
    "`r"

[void]$tr16.InsertAfter($newParas)

$full16 = $tr16.Text

# Bold "git push --follow-tags origin master"
$i = $full16.IndexOf("git push --follow-tags origin master")
$len = "git push --follow-tags origin master".Length
$tr16.Characters($i + 1, $len).Font.Bold = 1

# Bold only "pressing Ctrl + Shift + P" portion of the "In VS Code..." line
$i = $full16.IndexOf("pressing Ctrl + Shift + P")
$len = "pressing Ctrl + Shift + P".Length
$tr16.Characters($i + 1, $len).Font.Bold = 1

# Bold "Type Git: Push (Follow tags)"
$i = $full16.IndexOf("Type Git: Push (Follow tags)")
$len = "Type Git: Push (Follow tags)".Length
$tr16.Characters($i + 1, $len).Font.Bold = 1

# <a:ext cy="2893100"/> -> <a:ext cy="4416594"/>  (EMU -> points, 12700 EMU/pt)
# Must happen after the text edits above: spAutoFit recomputes the box
# height from the text content as soon as it changes, which would
# otherwise clobber an explicit size set earlier.
$tb16.Height = 347.7633570866142

# ---------------------------------------------------------------------------
# 2) "TextBox 19" (id=20) is now redundant - remove it entirely.
# ---------------------------------------------------------------------------
for ($idx = 1; $idx -le $s.Shapes.Count; $idx++) {
    $candidate = $s.Shapes.Item($idx)
    if ($candidate.Id -eq 20) {
        $candidate.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) "TextBox 20" (id=21) - bold the npm command and drop the extra space.
# ---------------------------------------------------------------------------
$tb20 = $null
for ($idx = 1; $idx -le $s.Shapes.Count; $idx++) {
    $candidate = $s.Shapes.Item($idx)
    if ($candidate.Id -eq 21) {
        $tb20 = $candidate
        break
    }
}

$tf20 = $tb20.TextFrame
$tr20 = $tf20.TextRange
$tr20.Text = "Use npm run build:dataload:watch to watch changes while developing."

$full20 = $tr20.Text

$i = $full20.IndexOf("npm run build:dataload:watch ")
$len = "npm run build:dataload:watch ".Length
$tr20.Characters($i + 1, $len).Font.Bold = 1
